$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Records to Create")
$ws2 = $wb.Worksheets.Item("Vocabularies")

# Rename some vocabulary / label strings before the column delete so the
# text substitutions are simpler to reason about (order doesn't matter for
# cells outside column D).

# 1. Delete column D ("Create DRs on Parent?") - this shifts everything
#    from column E onward one column to the left.
$ws1.Columns("D").Delete()

# 2. Update the row 1 "required/optional" banner cells (now shifted left).
$ws1.Range("B1").Value = "Optional"
$ws1.Range("E1:Q1").Value = "Optional (ITM, DO)"
$ws1.Range("G1").Value = "Required (ITM, DO)"
$ws1.Range("I1").Value = "Required (ITM, DO)"
$ws1.Range("R1:S1").Value = "Required (DO, DR_ATTACH)"

# 3. Update row 2 header for the record-type column.
$ws1.Range("B2").Value = "Record Type to Create (ITM, DO, DR_ATTACH)"

# 4. Update the record type values in the data rows.
$ws1.Range("B3:B38").Value = "DO"

# 5. Resize column B to fit the new, longer header text (COM ColumnWidth is
#    specified in "characters"; the stored XML width is chars + 5/6).
$ws1.Columns("B").ColumnWidth = 59.833333333333336

# 6. Reset sheet1's view back to the top-left / default selection.
$ws1.Range("A1").Select()

# 7. Update the Vocabularies "RecordType" list: DIGITAL_OBJECT -> DO, and
#    add the new DR_ATTACH entry below it (extends the RecordType table).
$ws2.Range("A3").Value = "DO"
$ws2.Range("A4").Value = "DR_ATTACH"
$ws2.ListObjects.Item("RecordType").Resize($ws2.Range("A2:A4"))

# 8. Update Vocabularies sheet selection to match the new active cell.
$ws2.Range("A5").Select()
